# Added Example for Sandbox apic
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the second L3OUT-TEST entry (row 4) with the new sandbox apic
# secondary-IP addresses and bump the VLAN numbers (1-4 -> 5-8).
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = "10.8.75.98"
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = "10.8.75.102"
$ws.Range("R4").Value = 7
$ws.Range("S4").Value = "10.8.75.106"
$ws.Range("V4").Value = 8
$ws.Range("W4").Value = "10.8.75.110"

# Widen column H slightly so the new data fits, matching the author's
# manual column resize (~13.33 characters).
$ws.Columns.Item(8).ColumnWidth = 12.5

# Scroll the view so column I is the left-most visible column, and move
# the active selection to V5 (just below the last data row).
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("V5").Select()
